$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Helper: set a plain text value, preserving percent-like strings ("20%")
# as literal text instead of letting Excel auto-convert them to a
# percentage number/format.
function Set-PercentText($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = $origStyle
}

# Row 2 (Mercedes)
$ws.Range("C2").Value = ("`$32.3m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.0m")
$ws.Range("D2").Value = "260 pts"

# Row 3 (Ferrari)
$ws.Range("C3").Value = ("`$26.2m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.1m")
$ws.Range("D3").Value = "97 pts"

# Row 4 (Red Bull)
Set-PercentText "B4" "20%"
$ws.Range("C4").Value = ("`$24.2m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.0m")
$ws.Range("D4").Value = "177 pts"

# Row 5 (McLaren)
$ws.Range("D5").Value = "108 pts"

# Row 6 (AlphaTauri)
$ws.Range("C6").Value = ("`$12.6m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.0m")
$ws.Range("D6").Value = "66 pts"

# Row 7 (Renault)
$ws.Range("C7").Value = ("`$12.1m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.0m")
$ws.Range("D7").Value = "103 pts"

# Row 8 (Racing Point)
Set-PercentText "B8" "18%"
$ws.Range("C8").Value = ("`$10.9m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.1m")
$ws.Range("D8").Value = "101 pts"

# Row 9 (Alfa Romeo)
$ws.Range("D9").Value = "62 pts"

# Row 10 (Haas)
$ws.Range("C10").Value = ("`$7.4m" + $nl + $nl + "      " + $nl + $nl + $nl + $nl + "    `$0.0m")
$ws.Range("D10").Value = "48 pts"

# Row 11 (Williams)
$ws.Range("D11").Value = "33 pts"
